$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fix the accent on "Michoacan" -> "Michoacán" for the three rows that
# reference that state (rows 11-13, column A).
$ws.Range("A11").Value = "Michoacán"
$ws.Range("A12").Value = "Michoacán"
$ws.Range("A13").Value = "Michoacán"

# Update the view: scroll so A7 is the top-left visible cell, and select
# the range of cells that were just edited (A11:A13), with A11 active.
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("A11:A13").Select()
